$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "59.443.20"
$ws.Cells.Item(2, 5).Value = "  +3.30%  "

$ws.Cells.Item(3, 4).Value = "2.395.69"
$ws.Cells.Item(3, 5).Value = "  +3.26%  "

$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  -0.12%  "

$ws.Cells.Item(5, 4).Value = "'552.87"
$ws.Cells.Item(5, 5).Value = "  +3.12%  "

$ws.Cells.Item(6, 4).Value = "'135.97"
$ws.Cells.Item(6, 5).Value = "  +1.56%  "

$ws.Cells.Item(7, 4).Value = "'1.00"
$ws.Cells.Item(7, 5).Value = "  -0.12%  "

$ws.Cells.Item(8, 5).Value = "  +2.78%  "

$ws.Cells.Item(10, 4).Value = "'5.85"
$ws.Cells.Item(10, 5).Value = "  +7.42%  "

$ws.Cells.Item(11, 5).Value = "  -0.99%  "

$ws.Cells.Item(12, 4).Value = "'0.362"
$ws.Cells.Item(12, 5).Value = "  +1.85%  "

$ws.Cells.Item(13, 4).Value = "'24.61"
$ws.Cells.Item(13, 5).Value = "  +4.71%  "

$ws.Cells.Item(14, 5).Value = "  +3.17%  "

$ws.Cells.Item(15, 4).Value = "59.348.90"
$ws.Cells.Item(15, 5).Value = "  +3.18%  "

$ws.Cells.Item(16, 5).Value = "  +5.46%  "

$ws.Cells.Item(17, 4).Value = "2.400.11"
$ws.Cells.Item(17, 5).Value = "  +3.06%  "

$ws.Cells.Item(18, 4).Value = "'11.25"
$ws.Cells.Item(18, 5).Value = "  +6.49%  "

$ws.Cells.Item(19, 4).Value = "'4.39"
$ws.Cells.Item(19, 5).Value = "  +4.19%  "

$ws.Cells.Item(20, 4).Value = "'335.72"
$ws.Cells.Item(20, 5).Value = "  +1.30%  "

$ws.Cells.Item(21, 4).Value = "'7.02"
$ws.Cells.Item(21, 5).Value = "  +5.97%  "

$ws.Cells.Item(22, 5).Value = "  -0.06%  "

$ws.Cells.Item(23, 4).Value = "'64.57"
$ws.Cells.Item(23, 5).Value = "  +3.97%  "

$ws.Cells.Item(24, 5).Value = "  +1.30%  "

$ws.Cells.Item(25, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(25, 4).Value = "'0.997"
$ws.Cells.Item(25, 5).Value = "  -0.29%  "

$ws.Cells.Item(26, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(26, 4).Value = "'8.44"
$ws.Cells.Item(26, 5).Value = "  +0.43%  "

$ws.Cells.Item(27, 4).Value = "'1.37"
$ws.Cells.Item(27, 5).Value = "  -1.59%  "

$ws.Cells.Item(28, 5).Value = "  +2.74%  "

$ws.Cells.Item(29, 4).Value = "0.0₃0763"
$ws.Cells.Item(29, 5).Value = "  +5.35%  "

$ws.Cells.Item(30, 4).Value = "'171.43"
$ws.Cells.Item(30, 5).Value = "  +0.88%  "

$ws.Cells.Item(31, 4).Value = "'6.25"
$ws.Cells.Item(31, 5).Value = "  +2.95%  "

$ws.Cells.Item(32, 4).Value = "'18.71"
$ws.Cells.Item(32, 5).Value = "  +1.87%  "

$ws.Cells.Item(33, 5).Value = "  -0.40%  "

$ws.Cells.Item(34, 4).Value = "'0.999"
$ws.Cells.Item(34, 5).Value = "  -0.04%  "

$ws.Cells.Item(35, 4).Value = "'4.30"
$ws.Cells.Item(35, 5).Value = "  +3.87%  "

$ws.Cells.Item(36, 4).Value = "'0.999"
$ws.Cells.Item(36, 5).Value = "  -0.10%  "

$ws.Cells.Item(37, 4).Value = "'1.28"
$ws.Cells.Item(37, 5).Value = "  +3.36%  "

$ws.Cells.Item(38, 4).Value = "'1.65"
$ws.Cells.Item(38, 5).Value = "  +2.69%  "

$ws.Cells.Item(39, 4).Value = "'40.13"
$ws.Cells.Item(39, 5).Value = "  +2.79%  "

$ws.Cells.Item(40, 4).Value = "'0.421"
$ws.Cells.Item(40, 5).Value = "  +12.53%  "

$ws.Cells.Item(41, 5).Value = "  +4.10%  "

$ws.Cells.Item(42, 4).Value = "'294.70"
$ws.Cells.Item(42, 5).Value = "  +4.56%  "

$ws.Cells.Item(43, 4).Value = "'141.90"
$ws.Cells.Item(43, 5).Value = "  -1.31%  "

$ws.Cells.Item(44, 4).Value = "'0.0965"
$ws.Cells.Item(44, 5).Value = "  +3.33%  "

$ws.Cells.Item(45, 4).Value = "'0.0525"
$ws.Cells.Item(45, 5).Value = "  +5.13%  "

$ws.Cells.Item(46, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(46, 4).Value = "'19.15"
$ws.Cells.Item(46, 5).Value = "  +1.42%  "

$ws.Cells.Item(47, 2).Value = "Mantle"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(47, 4).Value = "'0.570"
$ws.Cells.Item(47, 5).Value = "  +2.43%  "

$ws.Cells.Item(48, 4).Value = "'0.0226"
$ws.Cells.Item(48, 5).Value = "  +5.32%  "

$ws.Cells.Item(49, 4).Value = "'0.407"
$ws.Cells.Item(49, 5).Value = "  +6.95%  "

$ws.Cells.Item(50, 4).Value = "'11.03"
$ws.Cells.Item(50, 5).Value = "  -0.28%  "

$ws.Cells.Item(51, 4).Value = "'1.60"
$ws.Cells.Item(51, 5).Value = "  +5.36%  "

